# Updated envt to sandbox
# Replace the "neo*"/Michigan test-environment data in row 2 with the new
# "sushant*"/sandbox equivalents, update the main URL + cabi env values,
# fix up the hyperlink on A2 (drop its stale display text) and move the
# visible selection the way the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 data updates -----------------------------------------------
$ws.Range("A2").Value  = "https://sandbox.cabiclio.com/backoffice/control/main"
$ws.Range("E2").Value  = "sushanthost abc"
$ws.Range("F2").Value  = "sushantcohost abc"
$ws.Range("G2").Value  = "sushantguest1 abc"
$ws.Range("H2").Value  = "sushantguest2 abc"
$ws.Range("I2").Value  = "abc, sushanthost"
$ws.Range("J2").Value  = "abc, sushantcohost"
$ws.Range("L2").Value  = "abc, sushantguest1"
$ws.Range("M2").Value  = "abc, sushantguest2"
$ws.Range("N2").Value  = "cabisandbox"
$ws.Range("O2").Value  = 8

# --- Hyperlink on A2: keep the same target, drop the stale display text
$a2 = $ws.Range("A2")
$a2.Hyperlinks.Delete()
$ws.Hyperlinks.Add($a2, "https://test19.cliotest.com/backoffice/control/main") | Out-Null

# --- View state: scroll so column I is left-most and select N7 --------
$win = $excel.ActiveWindow
$win.ScrollColumn = 9
$win.ScrollRow = 1
$ws.Range("N7").Select()
